# Weekly update: insert a new data row (row 6) for the latest week,
# pushing all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 6 (shifts rows 6:47 -> 7:48).
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with this week's record.
$ws.Cells.Item(6, 1).Value  = 7
$ws.Cells.Item(6, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6, 3).Value  = "Ñuble"
$ws.Cells.Item(6, 4).Value  = 45190
$ws.Cells.Item(6, 5).Value  = 16
$ws.Cells.Item(6, 6).Value  = 300000000
$ws.Cells.Item(6, 7).Value  = "Espárragos"
$ws.Cells.Item(6, 8).Value  = "Sin especificar"
$ws.Cells.Item(6, 9).Value  = "Primera"
$ws.Cells.Item(6, 10).Value = 500
$ws.Cells.Item(6, 11).Value = 2000
$ws.Cells.Item(6, 12).Value = 2000
$ws.Cells.Item(6, 13).Value = 2000
$ws.Cells.Item(6, 14).Value = "$/kilo"
$ws.Cells.Item(6, 15).Value = "Región del Maule"
$ws.Cells.Item(6, 16).Value = 2000
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"
